# Create account api and graduate petition
# Adds a new generated account row (auto-username / shared password) to the
# Grade12StudentCredentials sheet, mirroring the pattern already used on the
# other grade-credential sheets in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade12StudentCredentials")

$ws.Cells.Item(2, 1).Value = "Auto20210827004426713"
$ws.Cells.Item(2, 2).Value = "Password@123"
